$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.013065630978526
$ws.Range("D2").Value = 1.053041699913072
$ws.Range("E2").Value = 1.015035990613206
$ws.Range("F2").Value = 1.054275782223558
$ws.Range("I2").Value = 1.04300669601295
$ws.Range("J2").Value = 1.018304738957197
$ws.Range("K2").Value = 1.055788763495031
$ws.Range("L2").Value = 1.017892351990943
$ws.Range("M2").Value = 1.057019440020605
$ws.Range("N2").Value = 1.010241486865392
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.013995088309817
$ws.Range("D3").Value = 1.053732133425986
$ws.Range("E3").Value = 1.015822506905327
$ws.Range("F3").Value = 1.055148175821361
$ws.Range("I3").Value = 1.043186372989367
$ws.Range("J3").Value = 1.018868276429063
$ws.Range("K3").Value = 1.056292477021201
$ws.Range("L3").Value = 1.018483995881102
$ws.Range("M3").Value = 1.057704893271765
$ws.Range("N3").Value = 1.010427958674937
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.014596390955547
$ws.Range("D4").Value = 1.054175760094066
$ws.Range("E4").Value = 1.016331751332403
$ws.Range("F4").Value = 1.055709279914853
$ws.Range("I4").Value = 1.043299420401619
$ws.Range("J4").Value = 1.019232229856665
$ws.Range("K4").Value = 1.056614612759628
$ws.Range("L4").Value = 1.018866494194132
$ws.Range("M4").Value = 1.058144402814282
$ws.Range("N4").Value = 1.010548361237007
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.014849149469311
$ws.Range("D5").Value = 1.054361509360631
$ws.Range("E5").Value = 1.016545912063937
$ws.Range("F5").Value = 1.055944353807938
$ws.Range("I5").Value = 1.043346174473119
$ws.Range("J5").Value = 1.019385069070498
$ws.Range("K5").Value = 1.056749126796056
$ws.Range("L5").Value = 1.019027215373162
$ws.Range("M5").Value = 1.058328207021805
$ws.Range("N5").Value = 1.010598916590716
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.014891587005504
$ws.Range("D6").Value = 1.054392653356732
$ws.Range("E6").Value = 1.016581874918337
$ws.Range("F6").Value = 1.05598377593982
$ws.Range("I6").Value = 1.043353979435117
$ws.Range("J6").Value = 1.019410721638663
$ws.Range("K6").Value = 1.0567716587661
$ws.Range("L6").Value = 1.019054196381867
$ws.Range("M6").Value = 1.058359011869985
$ws.Range("N6").Value = 1.010607401418814
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.014599768442618
$ws.Range("D7").Value = 1.0541782450401
$ws.Range("E7").Value = 1.016334612667933
$ws.Range("F7").Value = 1.055712424186019
$ws.Range("I7").Value = 1.043300048163242
$ws.Range("J7").Value = 1.01923427275835
$ws.Range("K7").Value = 1.056616413728312
$ws.Range("L7").Value = 1.018868642078351
$ws.Range("M7").Value = 1.058146862612199
$ws.Range("N7").Value = 1.010549037003693
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.013379769492723
$ws.Range("D8").Value = 1.053275682020736
$ws.Range("E8").Value = 1.015301731358788
$ws.Range("F8").Value = 1.05457131307287
$ws.Range("I8").Value = 1.043068083771302
$ws.Range("J8").Value = 1.018495332040073
$ws.Range("K8").Value = 1.055959781625164
$ws.Range("L8").Value = 1.018092369582117
$ws.Range("M8").Value = 1.057251924219556
$ws.Range("N8").Value = 1.010304558946729
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.011229101277119
$ws.Range("D9").Value = 1.051661365384887
$ws.Range("E9").Value = 1.013484134010088
$ws.Range("F9").Value = 1.052534627046993
$ws.Range("I9").Value = 1.042634757264553
$ws.Range("J9").Value = 1.017187954882145
$ws.Range("K9").Value = 1.054773686135148
$ws.Range("L9").Value = 1.016721954358397
$ws.Range("M9").Value = 1.055644193659167
$ws.Range("N9").Value = 1.009871801300646
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.009794774774515
$ws.Range("D10").Value = 1.050569204735547
$ws.Range("E10").Value = 1.01227413273856
$ws.Range("F10").Value = 1.051159523918303
$ws.Range("I10").Value = 1.042329415744889
$ws.Range("J10").Value = 1.01631287613076
$ws.Range("K10").Value = 1.053963544299537
$ws.Range("L10").Value = 1.015806699413919
$ws.Range("M10").Value = 1.054551819531851
$ws.Range("N10").Value = 1.009581998331026
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.009173572910528
$ws.Range("D11").Value = 1.050092525037843
$ws.Range("E11").Value = 1.011750612894063
$ws.Range("F11").Value = 1.050560001770569
$ws.Range("I11").Value = 1.042193308761241
$ws.Range("J11").Value = 1.015933138024868
$ws.Range("K11").Value = 1.053608159604286
$ws.Range("L11").Value = 1.015410003997909
$ws.Range("M11").Value = 1.054073955134469
$ws.Range("N11").Value = 1.009456205963612
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.00894281195868
$ws.Range("D12").Value = 1.049914900517055
$ws.Range("E12").Value = 1.0115562182483
$ws.Range("F12").Value = 1.050336698828908
$ws.Range("I12").Value = 1.042142168715633
$ws.Range("J12").Value = 1.015791963370464
$ws.Range("K12").Value = 1.053475466012004
$ws.Range("L12").Value = 1.015262596714774
$ws.Range("M12").Value = 1.053895726402764
$ws.Range("N12").Value = 1.0094094353806
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.00899231181517
$ws.Range("D13").Value = 1.049953027099592
$ws.Range("E13").Value = 1.011597913650234
$ws.Range("F13").Value = 1.05038462582072
$ws.Range("I13").Value = 1.0421531648357
$ws.Range("J13").Value = 1.015822251369246
$ws.Range("K13").Value = 1.053503960345017
$ws.Range("L13").Value = 1.015294218642409
$ws.Range("M13").Value = 1.053933990020323
$ws.Range("N13").Value = 1.009419469893505
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.009154498521633
$ws.Range("D14").Value = 1.050077854046007
$ws.Range("E14").Value = 1.011734542861886
$ws.Range("F14").Value = 1.050541555992502
$ws.Range("I14").Value = 1.04218909341453
$ws.Range("J14").Value = 1.015921470989574
$ws.Range("K14").Value = 1.053597205135587
$ws.Range("L14").Value = 1.015397820423276
$ws.Range("M14").Value = 1.054059237547601
$ws.Range("N14").Value = 1.009452340823464
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.009254424610656
$ws.Range("D15").Value = 1.050154689294461
$ws.Range("E15").Value = 1.011818733130857
$ws.Range("F15").Value = 1.050638164562862
$ws.Range("I15").Value = 1.042211152857573
$ws.Range("J15").Value = 1.015982587191303
$ws.Range("K15").Value = 1.053654565216355
$ws.Range("L15").Value = 1.01546164537625
$ws.Range("M15").Value = 1.054136310174867
$ws.Range("N15").Value = 1.009472587640102
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.009835999218071
$ws.Range("D16").Value = 1.050600761160548
$ws.Range("E16").Value = 1.012308885936296
$ws.Range("F16").Value = 1.051199226095431
$ws.Range("I16").Value = 1.042338366792355
$ws.Range("J16").Value = 1.016338060811174
$ws.Range("K16").Value = 1.053987033469413
$ws.Range("L16").Value = 1.015833018756299
$ws.Range("M16").Value = 1.05458343156389
$ws.Range("N16").Value = 1.009590340334753
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.010200771395824
$ws.Range("D17").Value = 1.050879562814147
$ws.Range("E17").Value = 1.012616458671837
$ws.Range("F17").Value = 1.051550070171794
$ws.Range("I17").Value = 1.042417123353286
$ws.Range("J17").Value = 1.016560820097831
$ws.Range("K17").Value = 1.054194354708083
$ws.Range("L17").Value = 1.016065869196448
$ws.Range("M17").Value = 1.054862599087475
$ws.Range("N17").Value = 1.009664121799075
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.010413524356977
$ws.Range("D18").Value = 1.051041819552614
$ws.Range("E18").Value = 1.012795900974541
$ws.Range("F18").Value = 1.051754316550228
$ws.Range("I18").Value = 1.042462685135309
$ws.Range("J18").Value = 1.016690672347655
$ws.Range("K18").Value = 1.054314838920716
$ws.Range("L18").Value = 1.016201649820826
$ws.Range("M18").Value = 1.055024963849049
$ws.Range("N18").Value = 1.009707127760893
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.010486065487648
$ws.Range("D19").Value = 1.051097083196887
$ws.Range("E19").Value = 1.012857092971486
$ws.Range("F19").Value = 1.051823892242296
$ws.Range("I19").Value = 1.04247815679162
$ws.Range("J19").Value = 1.016734935108442
$ws.Range("K19").Value = 1.054355845780671
$ws.Range("L19").Value = 1.016247941189367
$ws.Range("M19").Value = 1.05508024645333
$ws.Range("N19").Value = 1.009721786668261
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01016163608051
$ws.Range("D20").Value = 1.050849687627723
$ws.Range("E20").Value = 1.012583454849073
$ws.Range("F20").Value = 1.051512468773431
$ws.Range("I20").Value = 1.042408712361735
$ws.Range("J20").Value = 1.016536928349868
$ws.Range("K20").Value = 1.054172156887329
$ws.Range("L20").Value = 1.01604089037957
$ws.Range("M20").Value = 1.054832695547238
$ws.Range("N20").Value = 1.00965620879761
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.009106739097662
$ws.Range("D21").Value = 1.050041111191771
$ws.Range("E21").Value = 1.011694307191386
$ws.Range("F21").Value = 1.050495360900299
$ws.Range("I21").Value = 1.042178529453797
$ws.Range("J21").Value = 1.015892256674922
$ws.Range("K21").Value = 1.053569765847234
$ws.Range("L21").Value = 1.015367313847026
$ws.Range("M21").Value = 1.054022375374976
$ws.Range("N21").Value = 1.009442662413433
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.008443374241028
$ws.Range("D22").Value = 1.049529461998626
$ws.Range("E22").Value = 1.011135635166361
$ws.Range("F22").Value = 1.049852314565923
$ws.Range("I22").Value = 1.04203042674578
$ws.Range("J22").Value = 1.015486214597545
$ws.Range("K22").Value = 1.053187039389641
$ws.Range("L22").Value = 1.01494348044766
$ws.Range("M22").Value = 1.053508680002628
$ws.Range("N22").Value = 1.009308133052136
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.008795046709919
$ws.Range("D23").Value = 1.04980100593561
$ws.Range("E23").Value = 1.011431762329661
$ws.Range("F23").Value = 1.050193541717652
$ws.Range("I23").Value = 1.042109258663453
$ws.Range("J23").Value = 1.015701532402228
$ws.Range("K23").Value = 1.053390306673207
$ws.Range("L23").Value = 1.015168193503035
$ws.Range("M23").Value = 1.053781398714005
$ws.Range("N23").Value = 1.009379474583471
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.010179319683393
$ws.Range("D24").Value = 1.050863188061659
$ws.Range("E24").Value = 1.012598367731359
$ws.Range("F24").Value = 1.051529460447138
$ws.Range("I24").Value = 1.042412514087187
$ws.Range("J24").Value = 1.016547724247886
$ws.Range("K24").Value = 1.054182188494928
$ws.Range("L24").Value = 1.01605217734633
$ws.Range("M24").Value = 1.054846209118573
$ws.Range("N24").Value = 1.009659784433558
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.011785199938692
$ws.Range("D25").Value = 1.052081523861484
$ws.Range("E25").Value = 1.01395372644933
$ws.Range("F25").Value = 1.053064217809813
$ws.Range("I25").Value = 1.042749687897841
$ws.Range("J25").Value = 1.017526562523082
$ws.Range("K25").Value = 1.0550837489584
$ws.Range("L25").Value = 1.017076533179131
$ws.Range("M25").Value = 1.056063461017491
$ws.Range("N25").Value = 1.009983909517986

Write-Output "Applied vm_pu.xlsx value updates for rows 2-25"
